$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.910463666666666
$ws.Range("H2").Value = 11.731391
$ws.Range("I2").Value = 0.02584512419166262
$ws.Range("J2").Value = 0.02584512419166263
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 9.210619666666666
$ws.Range("N2").Value = 27.631859
$ws.Range("O2").Value = 0.133636377806767
$ws.Range("P2").Value = 0.133636377806767
$ws.Range("Q2").Value = 36.01779355398544
$ws.Range("R2").Value = 324.160141985869
$ws.Range("S2").Value = 0.00345384878093984
$ws.Range("T2").Value = 0.00345384878093984

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.910463666666666
$ws.Range("H3").Value = 11.731391
$ws.Range("I3").Value = 0.02584512419166262
$ws.Range("J3").Value = 0.02584512419166263
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 43.97212233333332
$ws.Range("N3").Value = 131.916367
$ws.Range("O3").Value = 0.6379891218794987
$ws.Range("P3").Value = 0.6379891218794989
$ws.Range("Q3").Value = 171.9513867307219
$ws.Range("R3").Value = 1547.562480576497
$ws.Range("S3").Value = 0.01648890808790543
$ws.Range("T3").Value = 0.01648890808790543

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.910463666666666
$ws.Range("H4").Value = 11.731391
$ws.Range("I4").Value = 0.02584512419166262
$ws.Range("J4").Value = 0.02584512419166263
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 15.740255
$ws.Range("N4").Value = 47.220765
$ws.Range("O4").Value = 0.2283745003137342
$ws.Range("P4").Value = 0.2283745003137342
$ws.Range("Q4").Value = 61.55169528156833
$ws.Range("R4").Value = 553.9652575341149
$ws.Range("S4").Value = 0.005902367322817355
$ws.Range("T4").Value = 0.005902367322817357

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 125.0119883333333
$ws.Range("H5").Value = 375.035965
$ws.Range("I5").Value = 0.826232037766454
$ws.Range("J5").Value = 0.8262320377664542
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.210619666666666
$ws.Range("N5").Value = 27.631859
$ws.Range("O5").Value = 0.133636377806767
$ws.Range("P5").Value = 0.133636377806767
$ws.Range("Q5").Value = 1151.437878312104
$ws.Range("R5").Value = 10362.94090480893
$ws.Range("S5").Value = 0.1104146567550128
$ws.Range("T5").Value = 0.1104146567550128

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 125.0119883333333
$ws.Range("H6").Value = 375.035965
$ws.Range("I6").Value = 0.826232037766454
$ws.Range("J6").Value = 0.8262320377664542
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 43.97212233333332
$ws.Range("N6").Value = 131.916367
$ws.Range("O6").Value = 0.6379891218794987
$ws.Range("P6").Value = 0.6379891218794989
$ws.Range("Q6").Value = 5497.042444126572
$ws.Range("R6").Value = 49473.38199713915
$ws.Range("S6").Value = 0.5271270522433289
$ws.Range("T6").Value = 0.527127052243329

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 125.0119883333333
$ws.Range("H7").Value = 375.035965
$ws.Range("I7").Value = 0.826232037766454
$ws.Range("J7").Value = 0.8262320377664542
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 15.740255
$ws.Range("N7").Value = 47.220765
$ws.Range("O7").Value = 0.2283745003137342
$ws.Range("P7").Value = 0.2283745003137342
$ws.Range("Q7").Value = 1967.720574423692
$ws.Range("R7").Value = 17709.48516981323
$ws.Range("S7").Value = 0.1886903287681123
$ws.Range("T7").Value = 0.1886903287681123

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 22.381277
$ws.Range("H8").Value = 67.14383099999999
$ws.Range("I8").Value = 0.1479228380418832
$ws.Range("J8").Value = 0.1479228380418833
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 9.210619666666666
$ws.Range("N8").Value = 27.631859
$ws.Range("O8").Value = 0.133636377806767
$ws.Range("P8").Value = 0.133636377806767
$ws.Range("Q8").Value = 206.1454301013143
$ws.Range("R8").Value = 1855.308870911829
$ws.Range("S8").Value = 0.01976787227081431
$ws.Range("T8").Value = 0.01976787227081432

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 22.381277
$ws.Range("H9").Value = 67.14383099999999
$ws.Range("I9").Value = 0.1479228380418832
$ws.Range("J9").Value = 0.1479228380418833
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 43.97212233333332
$ws.Range("N9").Value = 131.916367
$ws.Range("O9").Value = 0.6379891218794987
$ws.Range("P9").Value = 0.6379891218794989
$ws.Range("Q9").Value = 984.1522502202193
$ws.Range("R9").Value = 8857.370251981974
$ws.Range("S9").Value = 0.0943731615482644
$ws.Range("T9").Value = 0.09437316154826444

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 22.381277
$ws.Range("H10").Value = 67.14383099999999
$ws.Range("I10").Value = 0.1479228380418832
$ws.Range("J10").Value = 0.1479228380418833
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 15.740255
$ws.Range("N10").Value = 47.220765
$ws.Range("O10").Value = 0.2283745003137342
$ws.Range("P10").Value = 0.2283745003137342
$ws.Range("Q10").Value = 352.2870072056349
$ws.Range("R10").Value = 3170.583064850714
$ws.Range("S10").Value = 0.03378180422280452
$ws.Range("T10").Value = 0.03378180422280453
